# Generate Report for Handoff
# Swap the "174a5250..." and "9e8b9b78..." file rows on each sheet so the
# file that just got a new handoff ("174a5250...") now shows on row 3 with
# status "Ready for handoff" (and an updated "Latest Handoff Datetime"),
# while "9e8b9b78..." keeps its prior "Handed back" status and now sits on
# row 2.

$wb = $excel.ActiveWorkbook

$file174 = "174a5250-aaa8-4c7e-bb4b-37d715947d46.md"
$file9e8 = "9e8b9b78-eb94-40d6-bbd0-3b23b4f28699.md"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = $file9e8
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("C2").Value = "Handed back: in sync with en-US"

$ws.Range("A3").Value = $file174
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"

foreach ($h in $ws.Hyperlinks) {
    if ($h.Range.Row -eq 2 -and $h.Range.Column -eq 1) {
        $h.TextToDisplay = $file9e8
    } elseif ($h.Range.Row -eq 3 -and $h.Range.Column -eq 1) {
        $h.TextToDisplay = $file174
    }
}

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = $file9e8
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("C2").Value = "9e8b9b78-eb94-40d6-bbd0-3b23b4f28699.0499f178131fe92670f127d9da10a83ec3a86c70.zh-cn.xlf"
$ws.Range("D2").Value = "2016-03-09 08:19:50"
$ws.Range("E2").Value = $file9e8
$ws.Range("F2").Value = "9e8b9b78-eb94-40d6-bbd0-3b23b4f28699.0499f178131fe92670f127d9da10a83ec3a86c70.zh-cn.xlf"
$ws.Range("G2").Value = "2016-03-09 08:21:16"
$ws.Range("H2").Value = "Include"

$ws.Range("A3").Value = $file174
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "174a5250-aaa8-4c7e-bb4b-37d715947d46.ab6a47d2a271508164a1a32491a6817bf22610e4.zh-cn.xlf"
$ws.Range("D3").Value = "2016-03-09 08:21:56"
$ws.Range("E3").Value = $file174
$ws.Range("F3").Value = "174a5250-aaa8-4c7e-bb4b-37d715947d46.ab6a47d2a271508164a1a32491a6817bf22610e4.zh-cn.xlf"
$ws.Range("G3").Value = "2016-03-09 08:21:16"
$ws.Range("H3").Value = "Include"

foreach ($h in $ws.Hyperlinks) {
    $r = $h.Range.Row
    $c = $h.Range.Column
    if ($r -eq 2) {
        if ($c -eq 1 -or $c -eq 5) {
            $h.TextToDisplay = $file9e8
        } elseif ($c -eq 3 -or $c -eq 6) {
            $h.TextToDisplay = "9e8b9b78-eb94-40d6-bbd0-3b23b4f28699.0499f178131fe92670f127d9da10a83ec3a86c70.zh-cn.xlf"
        }
    } elseif ($r -eq 3) {
        if ($c -eq 1 -or $c -eq 5) {
            $h.TextToDisplay = $file174
        } elseif ($c -eq 3 -or $c -eq 6) {
            $h.TextToDisplay = "174a5250-aaa8-4c7e-bb4b-37d715947d46.ab6a47d2a271508164a1a32491a6817bf22610e4.zh-cn.xlf"
        }
    }
}

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = $file9e8
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("C2").Value = "9e8b9b78-eb94-40d6-bbd0-3b23b4f28699.0499f178131fe92670f127d9da10a83ec3a86c70.de-de.xlf"
$ws.Range("D2").Value = "2016-03-09 08:20:01"
$ws.Range("E2").Value = $file9e8
$ws.Range("F2").Value = "9e8b9b78-eb94-40d6-bbd0-3b23b4f28699.0499f178131fe92670f127d9da10a83ec3a86c70.de-de.xlf"
$ws.Range("G2").Value = "2016-03-09 08:21:27"
$ws.Range("H2").Value = "Include"

$ws.Range("A3").Value = $file174
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "174a5250-aaa8-4c7e-bb4b-37d715947d46.ab6a47d2a271508164a1a32491a6817bf22610e4.de-de.xlf"
$ws.Range("D3").Value = "2016-03-09 08:22:01"
$ws.Range("E3").Value = $file174
$ws.Range("F3").Value = "174a5250-aaa8-4c7e-bb4b-37d715947d46.ab6a47d2a271508164a1a32491a6817bf22610e4.de-de.xlf"
$ws.Range("G3").Value = "2016-03-09 08:21:27"
$ws.Range("H3").Value = "Include"

foreach ($h in $ws.Hyperlinks) {
    $r = $h.Range.Row
    $c = $h.Range.Column
    if ($r -eq 2) {
        if ($c -eq 1 -or $c -eq 5) {
            $h.TextToDisplay = $file9e8
        } elseif ($c -eq 3 -or $c -eq 6) {
            $h.TextToDisplay = "9e8b9b78-eb94-40d6-bbd0-3b23b4f28699.0499f178131fe92670f127d9da10a83ec3a86c70.de-de.xlf"
        }
    } elseif ($r -eq 3) {
        if ($c -eq 1 -or $c -eq 5) {
            $h.TextToDisplay = $file174
        } elseif ($c -eq 3 -or $c -eq 6) {
            $h.TextToDisplay = "174a5250-aaa8-4c7e-bb4b-37d715947d46.ab6a47d2a271508164a1a32491a6817bf22610e4.de-de.xlf"
        }
    }
}
